$d = $word.ActiveDocument

# --- 1) First paragraph: change text + append red-colored addition ---
$p1 = $d.Paragraphs(1).Range
$p1.Text = "This is a Microsoft word document.  "

$end = $p1.End - 1  # before the paragraph mark
$r = $d.Range($end, $end)
$r.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r2 = $d.Range($end, $end + ("(This is a change " + [char]0x2013 + " Ve").Length)
$r2.Font.Color = 255

$end2 = $r2.End
$r3Text = "rsion for main branch"
$r3 = $d.Range($end2, $end2)
$r3.InsertAfter($r3Text)
$r3b = $d.Range($end2, $end2 + $r3Text.Length)
$r3b.Font.Color = 255

$end3 = $r3b.End
$r4Text = ")"
$r4 = $d.Range($end3, $end3)
$r4.InsertAfter($r4Text)
$r4b = $d.Range($end3, $end3 + $r4Text.Length)
$r4b.Font.Color = 255

# --- 2) Remove the trailing "ank God almighty, we are free at last." paragraph ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete()
